# Fruta / hortaliza, semanal
# Refresh the weekly "Granada" price records (rows 2-13): dates, variety,
# quality, volumes, prices, units of sale, origin and derived $/kg all move
# to the latest weekly source-system snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  DateSerial = 44285; K = "Wonderfull";      L = "Primera"; M = 8;   N = 280000; O = 300000; P = 290000; Q = "`$/bins (400 kilos)";     R = "Provincia del Elquí";  S = 725;  T = 400 }
    @{ Row = 3;  DateSerial = 44662; K = "Sin especificar"; L = "Primera"; M = 45;  N = 18000;  O = 18000;  P = 18000;  Q = "`$/caja 18 kilos granel"; R = "Provincia de Limarí";  S = 1000; T = 18  }
    @{ Row = 4;  DateSerial = 44662; K = "Sin especificar"; L = "Segunda"; M = 60;  N = 16000;  O = 16000;  P = 16000;  Q = "`$/caja 18 kilos granel"; R = "Provincia de Limarí";  S = 889;  T = 18  }
    @{ Row = 5;  DateSerial = 44320; K = "Wonderfull";      L = "Primera"; M = 12;  N = 250000; O = 260000; P = 255000; Q = "`$/bins (400 kilos)";     R = "Provincia de Limarí";  S = 638;  T = 400 }
    @{ Row = 6;  DateSerial = 44307; K = "Sin especificar"; L = "Primera"; M = 150; N = 16000;  O = 18000;  P = 17000;  Q = "`$/caja 15 kilos granel"; R = "Región de O'Higgins";  S = 1133; T = 15  }
    @{ Row = 7;  DateSerial = 44280; K = "Sin especificar"; L = "Primera"; M = 15;  N = 360000; O = 360000; P = 360000; Q = "`$/bins (450 kilos)";     R = "Provincia del Elquí";  S = 800;  T = 450 }
    @{ Row = 8;  DateSerial = 44334; K = "Wonderfull";      L = "Primera"; M = 16;  N = 240000; O = 250000; P = 245000; Q = "`$/bins (450 kilos)";     R = "Provincia de Limarí";  S = 544;  T = 450 }
    @{ Row = 9;  DateSerial = 44312; K = "Wonderfull";      L = "Primera"; M = 24;  N = 220000; O = 240000; P = 230000; Q = "`$/bins (400 kilos)";     R = "Región de O'Higgins";  S = 575;  T = 400 }
    @{ Row = 10; DateSerial = 44312; K = "Wonderfull";      L = "Primera"; M = 34;  N = 240000; O = 240000; P = 240000; Q = "`$/bins (450 kilos)";     R = "Provincia del Elquí";  S = 533;  T = 450 }
    @{ Row = 11; DateSerial = 44266; K = "Wonderfull";      L = "Segunda"; M = 120; N = 4800;   O = 4800;   P = 4800;   Q = "`$/bandeja 4 kilos";      R = "Provincia del Elquí";  S = 1200; T = 4   }
    @{ Row = 12; DateSerial = 44266; K = "Wonderfull";      L = "Tercera"; M = 80;  N = 4000;   O = 4000;   P = 4000;   Q = "`$/bandeja 4 kilos";      R = "Provincia del Elquí";  S = 1000; T = 4   }
    @{ Row = 13; DateSerial = 44721; K = "Wonderfull";      L = "Primera"; M = 7;   N = 300000; O = 300000; P = 300000; Q = "`$/bins (400 kilos)";     R = "Región Metropolitana"; S = 750;  T = 400 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.DateSerial
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
